$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal TEXT value into a cell while preserving whatever
# number format / style the cell currently has (avoids the auto "numeric
# string -> number" coercion that a plain .Value = "0" assignment triggers).
# We do this by writing a string-literal formula, then collapsing it back to
# a plain value with a Copy / PasteSpecial(values-only) round-trip.
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

# ---------------------------------------------------------------------------
# Header strings: Volume/Number and the reporting week date range.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/22/2024  Through  1/28/2024"

# ---------------------------------------------------------------------------
# Row 15 (Rape): F15 flips from the number 1 to the text "0"; H15 updates.
# ---------------------------------------------------------------------------
$ws.Range("C15").Copy($ws.Range("F15"))
Set-TextValue $ws.Range("F15") "0"
$ws.Range("H15").Value = -100

# ---------------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -21.428571428571
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -21.428571428571
$ws.Range("L16").Value = -35.294117647058
$ws.Range("M16").Value = -45

# ---------------------------------------------------------------------------
# Row 17 (Fel. Assault): D17/E17 flip from text to numbers.
# ---------------------------------------------------------------------------
$ws.Range("C17").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 2
$ws.Range("H17").Copy($ws.Range("E17"))
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -11.111111111111
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = -11.111111111111
$ws.Range("L17").Value = 33.333333333333
$ws.Range("M17").Value = 33.333333333333
$ws.Range("N17").Value = -46.666666666666

# ---------------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = -35.483870967741
$ws.Range("I18").Value = 20
$ws.Range("J18").Value = 31
$ws.Range("K18").Value = -35.483870967741
$ws.Range("L18").Value = -20
$ws.Range("M18").Value = -23.076923076923
$ws.Range("N18").Value = -64.912280701754

# ---------------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 78
$ws.Range("G19").Value = 108
$ws.Range("H19").Value = -27.777777777777
$ws.Range("I19").Value = 78
$ws.Range("J19").Value = 108
$ws.Range("K19").Value = -27.777777777777
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 4
$ws.Range("N19").Value = -59.375

# ---------------------------------------------------------------------------
# Row 20 (G.L.A.): C20/D20/E20 flip from text to numbers.
# ---------------------------------------------------------------------------
$ws.Range("F20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("F20").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 1
$ws.Range("H20").Copy($ws.Range("E20"))
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 3
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -95.522388059701

# ---------------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 33
$ws.Range("E21").Value = -8.333333333333
$ws.Range("F21").Value = 120
$ws.Range("G21").Value = 165
$ws.Range("H21").Value = -27.272727272727
$ws.Range("I21").Value = 120
$ws.Range("J21").Value = 165
$ws.Range("K21").Value = -27.272727272727
$ws.Range("L21").Value = -10.447761194029
$ws.Range("M21").Value = -6.25
$ws.Range("N21").Value = -70.660146699266

# ---------------------------------------------------------------------------
# Row 22 (Transit): C22 flips from the number 1 to the text "0".
# ---------------------------------------------------------------------------
$ws.Range("C23").Copy($ws.Range("C22"))
Set-TextValue $ws.Range("C22") "0"
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -58.333333333333

# ---------------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 47
$ws.Range("E24").Value = -42.553191489361
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 146
$ws.Range("H24").Value = -31.506849315068
$ws.Range("I24").Value = 100
$ws.Range("J24").Value = 146
$ws.Range("K24").Value = -31.506849315068
$ws.Range("L24").Value = -25.925925925925
$ws.Range("M24").Value = -4.761904761904

# ---------------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -9.677419354838
$ws.Range("I25").Value = 28
$ws.Range("J25").Value = 31
$ws.Range("K25").Value = -9.677419354838
$ws.Range("M25").Value = 211.111111111111

# ---------------------------------------------------------------------------
# Row 26 (UCR Rape*): F26 flips from the number 2 to the text "0".
# ---------------------------------------------------------------------------
$ws.Range("C26").Copy($ws.Range("F26"))
Set-TextValue $ws.Range("F26") "0"
$ws.Range("H26").Value = -100

# ---------------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------------
$ws.Range("D27").Value = 3
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -75
